$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.741.14'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '2.103.74'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.46%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '346.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5196'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4426'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '54.04'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09360'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.174'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.86'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.49%  '
$ws.Range('D13').Value = '2.109.67'
$ws.Range('E13').Value = '  -0.38%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.301'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.82%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.810'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.21%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '102.57'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001159'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.73%  '
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '21.15'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.06673'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.315'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.006'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('D23').Value = '29.785.82'
$ws.Range('E23').Value = '  -1.24%  '
$ws.Range('E24').Value = '  -0.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.315'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.35%  '
$ws.Range('D26').Value = '2.359.31'
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.535'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '162.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.68'
$ws.Range('D30').Style = 'Normal'
$ws.Range('B31').Value = 'ARBITRUM'
$ws.Range('C31').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.797'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +9.67%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.145'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.1054'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.218'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.944'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.340'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.69'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02586'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06764'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.7003'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.59'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.338'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.80%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.2230'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6808'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.63%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.44'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.96%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.353'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.634'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.00000000356'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.51%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.210'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.78%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.220'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '81.88'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.56%  '
